# Updated Argent (silver) prices in Excel
# Each target cell holds its number as literal text (e.g. "5,455"), not a
# real numeric value, so we prefix the new value with a leading apostrophe
# to force Excel to keep storing it as text instead of re-interpreting the
# comma / decimal point and converting the cell to a number.
$wb = $excel.ActiveWorkbook

$ws6 = $wb.Worksheets.Item("Silver Rear_side")
$ws6.Range("B16").Value = "'5,440"

$ws7 = $wb.Worksheets.Item("Silver Busbar front-side")
$ws7.Range("B16").Value = "'8,145"

$ws8 = $wb.Worksheets.Item("Silver finger front-side")
$ws8.Range("B16").Value = "'8,195"

$ws9 = $wb.Worksheets.Item("USD_CNY")
$ws9.Range("B16").Value = "'7.2456"
